$d = $word.ActiveDocument

$replacements = @(
    @("918×7=6426", "867×4=3468"),
    @("469×4=1876", "817×7=5719"),
    @("184×4=736",  "688×3=2064"),
    @("330×8=2640", "937×7=6559"),
    @("666×2=1332", "755×2=1510"),
    @("238×2=476",  "286×7=2002"),
    @("824×9=7416", "160×4=640"),
    @("232×3=696",  "376×6=2256"),
    @("575×9=5175", "178×9=1602"),
    @("205×8=1640", "441×7=3087"),
    @("524×8=4192", "445×3=1335"),
    @("473×9=4257", "535×2=1070"),
    @("742×4=2968", "738×4=2952"),
    @("360×5=1800", "121×8=968"),
    @("683×8=5464", "734×2=1468"),
    @("751×6=4506", "709×9=6381"),
    @("369×6=2214", "297×3=891"),
    @("303×8=2424", "691×6=4146"),
    @("423×5=2115", "235×8=1880"),
    @("454×6=2724", "622×3=1866"),
    @("924×2=1848", "417×9=3753"),
    @("469×9=4221", "976×4=3904"),
    @("399×2=798",  "500×6=3000"),
    @("913×3=2739", "357×2=714"),
    @("987×7=6909", "449×7=3143")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
